{"js": "// Update the date line and the 25 division problems in the table to the\n// new values per the commit's regenerated worksheet data.\nconst replacements = [\n  [\"2024-04-26 Friday\", \"2024-04-27 Saturday\"],\n  [\"512\u00f74=\", \"703\u00f77=\"],\n  [\"898\u00f77=\", \"866\u00f72=\"],\n  [\"393\u00f79=\", \"789\u00f73=\"],\n  [\"789\u00f72=\", \"677\u00f77=\"],\n  [\"706\u00f78=\", \"510\u00f72=\"],\n  [\"801\u00f75=\", \"524\u00f77=\"],\n  [\"527\u00f76=\", \"861\u00f78=\"],\n  [\"437\u00f72=\", \"837\u00f79=\"],\n  [\"977\u00f75=\", \"338\u00f79=\"],\n  [\"122\u00f74=\", \"818\u00f73=\"],\n  [\"595\u00f76=\", \"647\u00f73=\"],\n  [\"704\u00f73=\", \"257\u00f75=\"],\n  [\"489\u00f72=\", \"986\u00f72=\"],\n  [\"933\u00f76=\", \"492\u00f72=\"],\n  [\"390\u00f75=\", \"461\u00f77=\"],\n  [\"690\u00f76=\", \"772\u00f79=\"],\n  [\"105\u00f72=\", \"382\u00f74=\"],\n  [\"653\u00f73=\", \"323\u00f74=\"],\n  [\"267\u00f73=\", \"494\u00f77=\"],\n  [\"461\u00f79=\", \"647\u00f77=\"],\n  [\"658\u00f77=\", \"349\u00f77=\"],\n  [\"922\u00f74=\", \"186\u00f75=\"],\n  [\"756\u00f79=\", \"823\u00f74=\"],\n  [\"587\u00f74=\", \"477\u00f73=\"],\n  [\"614\u00f75=\", \"258\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division problems in the table to the\n# new values per the commit's regenerated worksheet data.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"2024-04-26 Friday\", \"2024-04-27 Saturday\"),\n    @(\"512\u00f74=\", \"703\u00f77=\"),\n    @(\"898\u00f77=\", \"866\u00f72=\"),\n    @(\"393\u00f79=\", \"789\u00f73=\"),\n    @(\"789\u00f72=\", \"677\u00f77=\"),\n    @(\"706\u00f78=\", \"510\u00f72=\"),\n    @(\"801\u00f75=\", \"524\u00f77=\"),\n    @(\"527\u00f76=\", \"861\u00f78=\"),\n    @(\"437\u00f72=\", \"837\u00f79=\"),\n    @(\"977\u00f75=\", \"338\u00f79=\"),\n    @(\"122\u00f74=\", \"818\u00f73=\"),\n    @(\"595\u00f76=\", \"647\u00f73=\"),\n    @(\"704\u00f73=\", \"257\u00f75=\"),\n    @(\"489\u00f72=\", \"986\u00f72=\"),\n    @(\"933\u00f76=\", \"492\u00f72=\"),\n    @(\"390\u00f75=\", \"461\u00f77=\"),\n    @(\"690\u00f76=\", \"772\u00f79=\"),\n    @(\"105\u00f72=\", \"382\u00f74=\"),\n    @(\"653\u00f73=\", \"323\u00f74=\"),\n    @(\"267\u00f73=\", \"494\u00f77=\"),\n    @(\"461\u00f79=\", \"647\u00f77=\"),\n    @(\"658\u00f77=\", \"349\u00f77=\"),\n    @(\"922\u00f74=\", \"186\u00f75=\"),\n    @(\"756\u00f79=\", \"823\u00f74=\"),\n    @(\"587\u00f74=\", \"477\u00f73=\"),\n    @(\"614\u00f75=\", \"258\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceAll)\n}\n"}
